$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value  = 1022
$wsExpo.Range("F5").Value  = 20
$wsExpo.Range("F6").Value  = 460
$wsExpo.Range("F7").Value  = 743
$wsExpo.Range("F10").Value = 38
$wsExpo.Range("F14").Value = 860
$wsExpo.Range("F16").Value = 2005
$wsExpo.Range("F17").Value = 495
$wsExpo.Range("F18").Value = 7668
$wsExpo.Range("F19").Value = 561
$wsExpo.Range("F21").Value = 61
$wsExpo.Range("F24").Value = 229

# Sheet "本地生活" (Local Life)
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 5547
$wsLocal.Range("F3").Value = 401

# Sheet "全部类型" (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 5547
$wsAll.Range("F4").Value  = 401
$wsAll.Range("F7").Value  = 1023
$wsAll.Range("F11").Value = 20
$wsAll.Range("F12").Value = 460
$wsAll.Range("F13").Value = 743
$wsAll.Range("F17").Value = 38
$wsAll.Range("F23").Value = 860
$wsAll.Range("F26").Value = 2005
$wsAll.Range("F27").Value = 495
$wsAll.Range("F28").Value = 7668
$wsAll.Range("F31").Value = 561
$wsAll.Range("F33").Value = 61
$wsAll.Range("F37").Value = 229
